$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value2 = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "30.411.94"
Set-TextValue "E2" "  +0.85%  "
Set-TextValue "D3" "2.000.20"
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "324.18"
Set-TextValue "E5" "  +0.87%  "
Set-TextValue "D6" "0.9994"
Set-TextValue "E6" "  -0.05%  "
Set-TextValue "D7" "0.5109"
Set-TextValue "E7" "  +0.78%  "
Set-TextValue "D8" "0.4135"
Set-TextValue "E8" "  +2.22%  "
Set-TextValue "D9" "0.08720"
Set-TextValue "E9" "  +5.16%  "
Set-TextValue "D10" "1.136"
Set-TextValue "E10" "  +1.99%  "
Set-TextValue "D11" "42.89"
Set-TextValue "E11" "  +2.42%  "
Set-TextValue "D12" "24.78"
Set-TextValue "E12" "  +3.84%  "
Set-TextValue "D13" "1.996.89"
Set-TextValue "E13" "  +3.94%  "
Set-TextValue "D14" "6.541"
Set-TextValue "E14" "  +1.69%  "
Set-TextValue "D15" "7.438"
Set-TextValue "E15" "  +1.20%  "
Set-TextValue "D16" "0.9985"
Set-TextValue "E16" "  -0.21%  "
Set-TextValue "D17" "94.14"
Set-TextValue "E17" "  +1.36%  "
Set-TextValue "E18" "  +1.22%  "
Set-TextValue "D19" "0.06502"
Set-TextValue "E19" "  +0.25%  "
Set-TextValue "D20" "18.94"
Set-TextValue "E20" "  +1.17%  "
Set-TextValue "E21" "  -0.08%  "
Set-TextValue "D22" "6.175"
Set-TextValue "E22" "  +2.99%  "
Set-TextValue "D23" "30.458.46"
Set-TextValue "E23" "  +0.89%  "
Set-TextValue "D24" "11.85"
Set-TextValue "E24" "  +4.47%  "
Set-TextValue "D25" "2.210"
Set-TextValue "E25" "  +0.81%  "
Set-TextValue "D26" "2.229.99"
Set-TextValue "E26" "  +4.25%  "
Set-TextValue "D27" "22.44"
Set-TextValue "E27" "  +0.84%  "
Set-TextValue "D28" "163.35"
Set-TextValue "E28" "  +1.32%  "
Set-TextValue "D29" "2.423"
Set-TextValue "E29" "  +1.14%  "
Set-TextValue "D30" "131.59"
Set-TextValue "E30" "  +1.46%  "
Set-TextValue "D31" "1.143"
Set-TextValue "E31" "  +0.00%  "
Set-TextValue "D32" "0.1052"
Set-TextValue "E32" "  +0.59%  "
Set-TextValue "D33" "6.080"
Set-TextValue "E33" "  +0.99%  "
Set-TextValue "D34" "3.835"
Set-TextValue "E34" "  +1.27%  "
Set-TextValue "D35" "1.334"
Set-TextValue "E35" "  +11.30%  "
Set-TextValue "D36" "0.02511"
Set-TextValue "E36" "  +2.18%  "
Set-TextValue "D37" "0.06601"
Set-TextValue "E37" "  +1.81%  "
Set-TextValue "D38" "5.370"
Set-TextValue "E38" "  -1.47%  "
Set-TextValue "B39" "Aptos"
Set-TextValue "C39" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D39" "12.21"
Set-TextValue "E39" "  +6.94%  "
Set-TextValue "B40" "Algorand"
Set-TextValue "C40" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D40" "0.2202"
Set-TextValue "E40" "  +1.36%  "
Set-TextValue "D41" "9.034"
Set-TextValue "E41" "  +1.80%  "
Set-TextValue "D42" "0.6624"
Set-TextValue "E42" "  +3.00%  "
Set-TextValue "D43" "1.233"
Set-TextValue "E43" "  +1.09%  "
Set-TextValue "D44" "13.73"
Set-TextValue "E44" "  +1.73%  "
Set-TextValue "E45" "  +2.26%  "
Set-TextValue "D46" "2.208"
Set-TextValue "E46" "  +1.48%  "
Set-TextValue "D47" "3.663"
Set-TextValue "E47" "  +0.32%  "
Set-TextValue "D48" "1.261"
Set-TextValue "E48" "  +3.28%  "
Set-TextValue "D49" "124.52"
Set-TextValue "E49" "  +0.22%  "
Set-TextValue "D50" "80.32"
Set-TextValue "E50" "  +1.04%  "
Set-TextValue "D51" "0.06897"
Set-TextValue "E51" "  +1.17%  "
